# 10Th - MB for single stock and added new group
#
# This report keeps a rolling history of MarketBeat rank snapshots: each
# new trading day gets its own pair of columns inserted just to the right
# of the "Company" column (A), pushing all the previously captured days
# further to the right. Two brand-new brokerages (Benchmark, Evercore ISI)
# are also appended as new rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B. This shifts the existing
# Jun_17/Jun_15/Jun_13/Jun_10 columns (B:E) to the right (E:H) and carries
# their formatting/highlighting along with them.
$ws.Range("B1:D1").EntireColumn.Insert()

# New header row for the freshly inserted columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Default every data row in the new columns to "UN" (unchanged rating),
# matching the rest of the sheet's convention.
$ws.Range("B2:D27").Value = "UN"

# JPMorgan Chase & Co. (row 17) picked up a fresh rating note on both of
# the newly added Jun_26 columns.
$ws.Range("C17").Value = "6/20/2018,Reiterates,Neutral,GBX 1,500"
$ws.Range("D17").Value = "6/20/2018,Reiterates,Neutral,GBX 1,500"

# Two newly-tracked brokerages, added as new rows with no prior history.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
